# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "28.367.93"
Set-TextValue "D3" "1.565.85"
Set-TextValue "E3" "  -0.13%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "210.89"
Set-TextValue "E5" "  -0.47%  "
Set-TextValue "E6" "  -0.67%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "44.46"
Set-TextValue "E8" "  -3.47%  "
Set-TextValue "D9" "23.55"
Set-TextValue "E9" "  -1.96%  "
Set-TextValue "E10" "  -1.37%  "
Set-TextValue "E11" "  -0.64%  "
Set-TextValue "E12" "  +0.90%  "
Set-TextValue "D13" "1.787.73"
Set-TextValue "E13" "  -0.22%  "
Set-TextValue "D14" "1.577.62"
Set-TextValue "E14" "  +0.69%  "
Set-TextValue "E15" "  -0.33%  "
Set-TextValue "D16" "28.351.21"
Set-TextValue "E16" "  -0.49%  "
Set-TextValue "E17" "  -1.63%  "
Set-TextValue "D18" "60.23"
Set-TextValue "E18" "  -3.25%  "
Set-TextValue "D19" "228.02"
Set-TextValue "E19" "  +0.16%  "
Set-TextValue "E20" "  +0.36%  "
Set-TextValue "E21" "  -1.99%  "
Set-TextValue "E22" "  +0.01%  "
Set-TextValue "D23" "3.95"
Set-TextValue "E23" "  +1.39%  "
Set-TextValue "E24" "  -2.25%  "
Set-TextValue "D25" "2.05"
Set-TextValue "E25" "  -1.76%  "
Set-TextValue "D26" "150.18"
Set-TextValue "E26" "  -0.59%  "
Set-TextValue "E27" "  -0.86%  "
Set-TextValue "E28" "  +0.36%  "
Set-TextValue "E29" "  -2.19%  "
Set-TextValue "E30" "  -0.01%  "
Set-TextValue "D31" "0.0476"
Set-TextValue "E31" "  +1.85%  "
Set-TextValue "E32" "  -4.20%  "
Set-TextValue "E33" "  -1.13%  "
Set-TextValue "D34" "3.07"
Set-TextValue "E34" "  -0.05%  "
Set-TextValue "D35" "1.385.56"
Set-TextValue "E35" "  -0.75%  "
Set-TextValue "E36" "  +1.79%  "
Set-TextValue "E37" "  -3.34%  "
Set-TextValue "E38" "  -0.43%  "
Set-TextValue "D39" "2.64"
Set-TextValue "E39" "  +2.17%  "
Set-TextValue "E40" "  -2.22%  "
Set-TextValue "D41" "1.94"
Set-TextValue "E41" "  +3.55%  "
Set-TextValue "D42" "0.517"
Set-TextValue "E42" "  -3.60%  "
Set-TextValue "E43" "  -0.01%  "
Set-TextValue "E44" "  -0.19%  "
Set-TextValue "E45" "  -1.50%  "
Set-TextValue "D46" "5.35"
Set-TextValue "E46" "  -2.93%  "
Set-TextValue "B47" "Aave"
Set-TextValue "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "62.19"
Set-TextValue "E47" "  -1.14%  "
Set-TextValue "B48" "WEMIXToken"
Set-TextValue "C48" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D48" "0.920"
Set-TextValue "E48" "  -5.51%  "
Set-TextValue "D49" "1.700.71"
Set-TextValue "D50" "85.35"
Set-TextValue "E50" "  -0.67%  "
Set-TextValue "D51" "0.0₆0101"
Set-TextValue "E51" "  -1.65%  "
